$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Integrity" header in G1, matching the formatting of the
# existing header cells (bold, centered, bordered) by copying the style
# from the neighboring F1 header cell.
$ws.Range("G1").Value = "Integrity"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

$lastRow = $ws.UsedRange.Rows.Count

# Populate the Integrity column with a count of how many of the
# C:F (wire/inverter reading) cells are populated for that row.
for ($r = 2; $r -le $lastRow; $r++) {
    $count = 0
    for ($col = 3; $col -le 6; $col++) {
        $val = $ws.Cells.Item($r, $col).Value()
        if (-not [string]::IsNullOrEmpty($val)) {
            $count = $count + 1
        }
    }
    $ws.Cells.Item($r, 7).Value = $count
}
